$wb = $excel.ActiveWorkbook

# ALC row 19: Unbreak My Heart | Roof Tile
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2687.6667
$ws.Range("I19").Value = 1371
$ws.Range("J19").Value = 4004.3333
$ws.Range("K19").Value = 1371
$ws.Range("L19").Value = 4004.3333
$ws.Range("M19").Value = -1196
$ws.Range("N19").Value = -4354.3333

# ALC row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4573.467
$ws.Range("I62").Value = 4605.909
$ws.Range("J62").Value = 4484.25
$ws.Range("K62").Value = 4605.909
$ws.Range("L62").Value = 4484.25
$ws.Range("M62").Value = -3981.909
$ws.Range("N62").Value = -5732.25

# ALC row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4573.467
$ws.Range("I65").Value = 4605.909
$ws.Range("J65").Value = 4484.25
$ws.Range("K65").Value = 23029.545
$ws.Range("L65").Value = 22421.25
$ws.Range("M65").Value = -19909.545
$ws.Range("N65").Value = -28661.25

# ALC row 116: Growing Up | Growth Formula Kappa
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6695
$ws.Range("I116").Value = 6504.1714
$ws.Range("J116").Value = 7087.8823
$ws.Range("K116").Value = 6504.1714
$ws.Range("L116").Value = 7087.8823
$ws.Range("M116").Value = -3062.1714
$ws.Range("N116").Value = -13971.8823

# ALC row 121: Mindful Medicine | Tincture of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 3660
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3660
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 10980
$ws.Range("N121").Value = -14474

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3311.92
$ws.Range("I132").Value = 3318.5264
$ws.Range("J132").Value = 3291
$ws.Range("K132").Value = 9955.5792
$ws.Range("L132").Value = 9873
$ws.Range("M132").Value = -7425.5792
$ws.Range("N132").Value = -14933

# ALC row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 154332.17
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 154332.17
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 154332.17
$ws.Range("N133").Value = -164452.17

# ALC row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 12195915
$ws.Range("I135").Value = 14706699
$ws.Range("J135").Value = 681.5714
$ws.Range("K135").Value = 132360291
$ws.Range("L135").Value = 6134.1426
$ws.Range("M135").Value = -132357756
$ws.Range("N135").Value = -11204.1426

# ALC row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1195
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1195
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 3585
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -8685

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4924.25
$ws.Range("I138").Value = 3100
$ws.Range("J138").Value = 5090.091
$ws.Range("K138").Value = 9300
$ws.Range("L138").Value = 15270.273
$ws.Range("M138").Value = -4160
$ws.Range("N138").Value = -25550.273

# ALC row 140: Tome for Tradition | Book of Ra'Kaznar
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 72533.42999999999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 72533.42999999999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 72533.42999999999
$ws.Range("N140").Value = -82893.42999999999

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1843.25
$ws.Range("I74").Value = 1790.2
$ws.Range("J74").Value = 1931.6666
$ws.Range("K74").Value = 1790.2
$ws.Range("L74").Value = 1931.6666
$ws.Range("M74").Value = -916.2
$ws.Range("N74").Value = -3679.6666

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1843.25
$ws.Range("I77").Value = 1790.2
$ws.Range("J77").Value = 1931.6666
$ws.Range("K77").Value = 8951
$ws.Range("L77").Value = 9658.333000000001
$ws.Range("M77").Value = -4583
$ws.Range("N77").Value = -18394.333

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8283.796
$ws.Range("I31").Value = 2580.2222
$ws.Range("J31").Value = 11595.549
$ws.Range("K31").Value = 2580.2222
$ws.Range("L31").Value = 11595.549
$ws.Range("M31").Value = -2285.2222

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8283.796
$ws.Range("I34").Value = 2580.2222
$ws.Range("J34").Value = 11595.549
$ws.Range("K34").Value = 2580.2222
$ws.Range("L34").Value = 11595.549
$ws.Range("M34").Value = -2378.2222

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4055.5
$ws.Range("I58").Value = 4115.0713
$ws.Range("J58").Value = 3847
$ws.Range("K58").Value = 4115.0713
$ws.Range("L58").Value = 3847
$ws.Range("M58").Value = -3912.0713

# CRP row 135: The Wing's Wings | Ceiba Wings
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 69311.11
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 69311.11
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 69311.11
$ws.Range("N135").Value = -79451.11

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4055.5
$ws.Range("I136").Value = 4115.0713
$ws.Range("J136").Value = 3847
$ws.Range("K136").Value = 12345.2139
$ws.Range("L136").Value = 11541
$ws.Range("M136").Value = -9795.213899999999

# CRP row 138: Bow Out | Acacia Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 78487.664
$ws.Range("I138").Value = 78000
$ws.Range("J138").Value = 78548.625
$ws.Range("K138").Value = 78000
$ws.Range("L138").Value = 78548.625
$ws.Range("M138").Value = -72860
$ws.Range("N138").Value = -88828.625

# CUL row 68: Such a Butter Face | Fermented Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3075.9656
$ws.Range("I68").Value = 2500.6
$ws.Range("J68").Value = 3195.8333
$ws.Range("K68").Value = 7501.799999999999
$ws.Range("L68").Value = 9587.499899999999
$ws.Range("M68").Value = -6690.799999999999
$ws.Range("N68").Value = -11209.4999

# CUL row 71: No Margarine of Error (L) | Fermented Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3075.9656
$ws.Range("I71").Value = 2500.6
$ws.Range("J71").Value = 3195.8333
$ws.Range("K71").Value = 22505.4
$ws.Range("L71").Value = 28762.4997
$ws.Range("M71").Value = -18449.4
$ws.Range("N71").Value = -36874.4997

# CUL row 137: Creative Chocolate | Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2060.4167
$ws.Range("I137").Value = 1241.3334
$ws.Range("J137").Value = 2879.5
$ws.Range("K137").Value = 3724.0002
$ws.Range("L137").Value = 8638.5
$ws.Range("M137").Value = 1375.9998

# GSM row 122: Awarding Academic Excellence | Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 61410.555
$ws.Range("I122").Value = 89143.13
$ws.Range("J122").Value = 12345.23
$ws.Range("K122").Value = 267429.39
$ws.Range("L122").Value = 37035.69
$ws.Range("M122").Value = -264979.39
$ws.Range("N122").Value = -41935.69

# GSM row 126: Gold Rush Order | Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3152.9
$ws.Range("I126").Value = 3055.9473
$ws.Range("J126").Value = 4995
$ws.Range("K126").Value = 9167.841899999999
$ws.Range("L126").Value = 14985
$ws.Range("M126").Value = -6697.841899999999

# GSM row 135: Fan of the Foreign | Ruthenium Folding Fans
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 97444.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 97444.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 97444.25
$ws.Range("N135").Value = -107584.25

# LTW row 100: Tiger in the Sack | Tiger Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2308.2058
$ws.Range("I100").Value = 1938.7778
$ws.Range("J100").Value = 2723.8125
$ws.Range("K100").Value = 1938.7778
$ws.Range("L100").Value = 2723.8125
$ws.Range("M100").Value = -1397.7778

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 30000
$ws.Range("I132").Value = 26250
$ws.Range("J132").Value = 35000
$ws.Range("K132").Value = 78750
$ws.Range("L132").Value = 105000
$ws.Range("M132").Value = -76220

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 12823550
$ws.Range("I136").Value = 2823.4375
$ws.Range("J136").Value = 33336712
$ws.Range("K136").Value = 8470.3125
$ws.Range("L136").Value = 100010136
$ws.Range("M136").Value = -5920.3125
$ws.Range("N136").Value = -100015236

# WVR row 113: A Tender Table | Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 776.6889
$ws.Range("I113").Value = 850.6667
$ws.Range("J113").Value = 480.77777
$ws.Range("K113").Value = 2552.0001
$ws.Range("L113").Value = 1442.33331
$ws.Range("M113").Value = -382.0001000000002
$ws.Range("N113").Value = -5782.33331

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 31667.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 31667.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 95002.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -100062.5

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2023.9546
$ws.Range("I136").Value = 1780.6428
$ws.Range("J136").Value = 2449.75
$ws.Range("K136").Value = 5341.928400000001
$ws.Range("L136").Value = 7349.25
$ws.Range("M136").Value = -2791.928400000001
$ws.Range("N136").Value = -12449.25

# WVR row 141: Silk for Sunperch | Thunderyards Silk Coat of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 102880
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 102880
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 102880
$ws.Range("N141").Value = -113240
